$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook tracks a rolling 10-quarter window (columns E:N).
# This update drops the oldest quarter (فصل دوم منتهی به 1399/06) and
# appends the newest quarter (فصل چهارم منتهی به 1401/12), shifting all
# quarterly figures one column to the left and refreshing the new-quarter values.

# --- Quarter header labels (rows 8, 18, 32, 47, 60, 72) ---
$quarterHeaders = @("فصل سوم منتهی به 1399/09", "فصل چهارم منتهی به 1399/12", "فصل اول منتهی به 1400/03", "فصل دوم منتهی به 1400/06", "فصل سوم منتهی به 1400/09", "فصل چهارم منتهی به 1400/12", "فصل اول منتهی به 1401/03", "فصل دوم منتهی به 1401/06", "فصل سوم منتهی به 1401/09", "فصل چهارم منتهی به 1401/12")
$headerRows = @(8, 18, 32, 47, 60, 72)
foreach ($r in $headerRows) {
    for ($i = 0; $i -lt $quarterHeaders.Length; $i++) {
        $ws.Cells.Item($r, 5 + $i).Value = $quarterHeaders[$i]
    }
}

# --- Data rows: new values for columns E:N ---
$rowData = @{
    11 = @("-", "-", "-", 0, 0, "-", 0, 0, 0, 0)
    12 = @("-", "-", "-", 15924011, -15906702, "-", "-", "-", "-", "-")
    13 = @(6059, 5365, 7641, "-", "-", 8062, 8135, 7252, 7292, 6140)
    14 = @(6059, 5365, 7641, 15924011, -15906702, 8062, 8135, 7252, 7292, 6140)
    20 = @(-13241419, 1806512, "-", "-", "-", "-", "-", "-", "-", "-")
    21 = @(5372681, 24936255, "-", "-", "-", "-", "-", "-", "-", "-")
    22 = @(0, 0, "-", "-", "-", "-", "-", "-", "-", "-")
    23 = @(0, 0, "-", 0, 0, "-", 0, 0, 0, 0)
    24 = @("-", "-", 7549380, 8515982, 8743, "-", "-", "-", "-", "-")
    25 = @("-", "-", "-", "-", "-", "-", "-", "-", "-", "-")
    26 = @("-", "-", "-", "-", "-", "-", "-", "-", "-", "-")
    27 = @("-", "-", "-", "-", "-", 8032, 8006, 6577, 7517, 6994)
    28 = @(-7868738, 26742767, 7549380, 8515982, 8743, 8032, 8006, 6577, 7517, 6994)
    35 = @(-1229541, -2117834, "-", "-", "-", "-", "-", "-", "-", "-")
    36 = @(788073, 748138, "-", "-", "-", "-", "-", "-", "-", "-")
    37 = @(0, 0, "-", "-", "-", "-", "-", "-", "-", "-")
    38 = @(0, 0, "-", 0, 0, "-", 0, 0, 0, 0)
    39 = @("-", "-", 1721437, 2056560, 2295244, "-", "-", "-", "-", "-")
    40 = @("-", "-", "-", "-", "-", "-", "-", "-", "-", "-")
    41 = @("-", "-", "-", "-", "-", "-", "-", "-", "-", "-")
    42 = @("-", "-", "-", "-", "-", 2236888, 3138299, 2787141, 2675641, 2430474)
    43 = @(-441468, -1369696, 1721437, 2056560, 2295244, 2236888, 3138299, 2787141, 2675641, 2430474)
    49 = @(134923, 208621, "-", "-", "-", "-", "-", "-", "-", "-")
    50 = @(145420, 189600, "-", "-", "-", "-", "-", "-", "-", "-")
    51 = @("-", "-", "-", "-", "-", "-", "-", "-", "-", "-")
    52 = @("-", "-", "-", "-", "-", "-", "-", "-", "-", "-")
    53 = @("-", "-", 228024, 241494, -143090, "-", "-", "-", "-", "-")
    54 = @("-", "-", "-", "-", "-", "-", "-", "-", "-", "-")
    55 = @("-", "-", "-", "-", "-", "-", "-", "-", "-", "-")
    56 = @("-", "-", "-", "-", "-", 253048995, 391993380, 423770868, 355945324, 347508436)
    62 = @(9336, -197898, "-", "-", "-", "-", "-", "-", "-", "-")
    63 = @(-747442, -532770, "-", "-", "-", "-", "-", "-", "-", "-")
    64 = @(0, 0, "-", "-", "-", "-", "-", "-", "-", "-")
    65 = @(0, 0, "-", 0, 0, "-", 0, 0, 0, 0)
    66 = @("-", "-", -1266016, -1548633, -1830721, "-", "-", "-", "-", "-")
    67 = @("-", "-", "-", "-", "-", -1931301, -2484291, -2170168, -2150970, -1976307)
    68 = @(-738106, -730668, -1266016, -1548633, -1830721, -1931301, -2484291, -2170168, -2150970, -1976307)
    74 = @(182210, -112236, "-", "-", "-", "-", "-", "-", "-", "-")
    75 = @(53639, 401250, "-", "-", "-", "-", "-", "-", "-", "-")
    76 = @(0, 0, "-", "-", "-", "-", "-", "-", "-", "-")
    77 = @(0, 0, "-", 0, 0, "-", 0, 0, 0, 0)
    78 = @("-", "-", 455421, 507927, 464523, "-", "-", "-", "-", "-")
    79 = @("-", "-", "-", "-", "-", 305587, 654008, 616973, 524671, 454167)
    80 = @(235849, 289014, 455421, 507927, 464523, 305587, 654008, 616973, 524671, 454167)
}
foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item([int]$r, 5 + $i).Value = $vals[$i]
    }
}
